# Update the "Numb" (column H) value for specific rows from 0 to -1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(31, 32, 34, 41, 43)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = -1
}
